$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2899
$ws.Range("I29").Value = 2898
$ws.Range("K29").Value = 8694
$ws.Range("M29").Value = -8413
$ws.Range("H51").Value = 6000
$ws.Range("J51").Value = 6000
$ws.Range("L51").Value = 6000
$ws.Range("N51").Value = -6968
$ws.Range("H53").Value = 138.57143
$ws.Range("J53").Value = 166
$ws.Range("L53").Value = 166
$ws.Range("N53").Value = -1440
$ws.Range("H88").Value = 7179.5
$ws.Range("I88").Value = 7997
$ws.Range("J88").Value = 7088.6665
$ws.Range("K88").Value = 7997
$ws.Range("L88").Value = 7088.6665
$ws.Range("M88").Value = -7591
$ws.Range("N88").Value = -7900.6665
$ws.Range("H91").Value = 7179.5
$ws.Range("I91").Value = 7997
$ws.Range("J91").Value = 7088.6665
$ws.Range("K91").Value = 7997
$ws.Range("L91").Value = 7088.6665
$ws.Range("M91").Value = -6593
$ws.Range("N91").Value = -9896.666499999999
$ws.Range("H98").Value = 757.5
$ws.Range("I98").Value = 607.6
$ws.Range("K98").Value = 607.6
$ws.Range("M98").Value = 890.4
$ws.Range("H112").Value = 2424
$ws.Range("I112").Value = 1995
$ws.Range("J112").Value = 2509.8
$ws.Range("K112").Value = 5985
$ws.Range("L112").Value = 7529.400000000001
$ws.Range("M112").Value = -4877
$ws.Range("N112").Value = -9745.400000000001
$ws.Range("H122").Value = 757.5
$ws.Range("I122").Value = 607.6
$ws.Range("K122").Value = 1822.8
$ws.Range("M122").Value = 627.1999999999998
$ws.Range("H135").Value = 2976.6
$ws.Range("I135").Value = 2996.5
$ws.Range("J135").Value = 2963.3333
$ws.Range("K135").Value = 26968.5
$ws.Range("L135").Value = 26669.9997
$ws.Range("M135").Value = -24433.5
$ws.Range("N135").Value = -31739.9997
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12460.9
$ws.Range("I32").Value = 8281.4
$ws.Range("K32").Value = 8281.4
$ws.Range("M32").Value = -7994.4
$ws.Range("H45").Value = 1452.75
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H110").Value = 9892.714
$ws.Range("I110").Value = 6330
$ws.Range("K110").Value = 6330
$ws.Range("M110").Value = -4285
$ws.Range("H132").Value = 2209.2
$ws.Range("I132").Value = 2209.2
$ws.Range("K132").Value = 6627.599999999999
$ws.Range("M132").Value = -4097.599999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4869.7144
$ws.Range("I86").Value = 4596.2
$ws.Range("J86").Value = 5553.5
$ws.Range("K86").Value = 4596.2
$ws.Range("L86").Value = 5553.5
$ws.Range("M86").Value = -3473.2
$ws.Range("N86").Value = -7799.5
$ws.Range("H89").Value = 4869.7144
$ws.Range("I89").Value = 4596.2
$ws.Range("J89").Value = 5553.5
$ws.Range("K89").Value = 22981
$ws.Range("L89").Value = 27767.5
$ws.Range("M89").Value = -17365
$ws.Range("N89").Value = -38999.5
$ws.Range("H107").Value = 683.75
$ws.Range("I107").Value = 683.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 683.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1236.25
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 2032.1482
$ws.Range("J134").Value = 3599.5
$ws.Range("L134").Value = 10798.5
$ws.Range("N134").Value = -15868.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 3974
$ws.Range("I36").Value = 3974
$ws.Range("K36").Value = 3974
$ws.Range("M36").Value = -3586
$ws.Range("H40").Value = 3974
$ws.Range("I40").Value = 3974
$ws.Range("K40").Value = 3974
$ws.Range("M40").Value = -3814
$ws.Range("H122").Value = 1978.4138
$ws.Range("I122").Value = 2012
$ws.Range("J122").Value = 1525
$ws.Range("K122").Value = 6036
$ws.Range("L122").Value = 4575
$ws.Range("M122").Value = -3586
$ws.Range("N122").Value = -9475
$ws.Range("H132").Value = 3803
$ws.Range("I132").Value = 3803
$ws.Range("K132").Value = 11409
$ws.Range("M132").Value = -8879
$ws.Range("H134").Value = 2961.0952
$ws.Range("I134").Value = 2315.6155
$ws.Range("K134").Value = 6946.8465
$ws.Range("M134").Value = -4411.8465
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H137").Value = 5644.273
$ws.Range("I137").Value = 3723.7144
$ws.Range("J137").Value = 9005.25
$ws.Range("K137").Value = 11171.1432
$ws.Range("L137").Value = 27015.75
$ws.Range("M137").Value = -6071.143199999999
$ws.Range("N137").Value = -37215.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5900
$ws.Range("J80").Value = 8750
$ws.Range("L80").Value = 8750
$ws.Range("N80").Value = -10746
$ws.Range("H83").Value = 5900
$ws.Range("J83").Value = 8750
$ws.Range("L83").Value = 43750
$ws.Range("N83").Value = -53734
$ws.Range("H107").Value = 643.2727
$ws.Range("I107").Value = 354
$ws.Range("J107").Value = 808.5714
$ws.Range("K107").Value = 354
$ws.Range("L107").Value = 808.5714
$ws.Range("M107").Value = 1566
$ws.Range("N107").Value = -4648.5714
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7374.625
$ws.Range("I22").Value = 1432.6666
$ws.Range("J22").Value = 10939.8
$ws.Range("K22").Value = 1432.6666
$ws.Range("L22").Value = 10939.8
$ws.Range("M22").Value = -1137.6666
$ws.Range("N22").Value = -11529.8
$ws.Range("H27").Value = 7374.625
$ws.Range("I27").Value = 1432.6666
$ws.Range("J27").Value = 10939.8
$ws.Range("K27").Value = 1432.6666
$ws.Range("L27").Value = 10939.8
$ws.Range("M27").Value = -1325.6666
$ws.Range("N27").Value = -11153.8
$ws.Range("H82").Value = 2534.7
$ws.Range("J82").Value = 1508.3636
$ws.Range("L82").Value = 1508.3636
$ws.Range("N82").Value = -2230.3636
$ws.Range("H85").Value = 2534.7
$ws.Range("J85").Value = 1508.3636
$ws.Range("L85").Value = 1508.3636
$ws.Range("N85").Value = -4004.3636
$ws.Range("H136").Value = 3927.8572
$ws.Range("I136").Value = 3599.2
$ws.Range("J136").Value = 4749.5
$ws.Range("K136").Value = 10797.6
$ws.Range("L136").Value = 14248.5
$ws.Range("M136").Value = -8247.599999999999
$ws.Range("N136").Value = -19348.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 18500
$ws.Range("J15").Value = 18500
$ws.Range("L15").Value = 18500
$ws.Range("N15").Value = -19076
$ws.Range("H81").Value = 1611
$ws.Range("I81").Value = 1611
$ws.Range("K81").Value = 3222
$ws.Range("M81").Value = -2161
$ws.Range("H84").Value = 1611
$ws.Range("I84").Value = 1611
$ws.Range("K84").Value = 16110
$ws.Range("M84").Value = -10806
$ws.Range("H136").Value = 1586.525
$ws.Range("I136").Value = 1207.4722
$ws.Range("K136").Value = 3622.4166
$ws.Range("M136").Value = -1072.4166
